# doc todo list updates
#
# 1. Window geometry tweak (workbookView xWindow/yWindow/windowWidth/windowHeight,
#    and dropping the stale activeTab since "asis" becomes the active sheet again).
# 2. "asis" sheet (sheet1): becomes the tab-selected sheet, the frozen pane scrolls
#    back to the top (A2, since row 1 is frozen), and the selection moves to B11.
# 3. "list" sheet (sheet3): a new blank column is inserted before column A, shifting
#    all existing data from columns A:C to B:D, and the selection moves to G12.

$wb = $excel.ActiveWorkbook

# --- window geometry -------------------------------------------------------
$win = $excel.ActiveWindow
$win.Left = 240
$win.Top = 210
$win.Width = 7215
$win.Height = 8070

# --- "list" sheet: shift data one column to the right -----------------------
$wsList = $wb.Worksheets.Item("list")
$wsList.Columns("A:A").Insert()
$wsList.Range("G12").Select()

# --- "asis" sheet: make it the active tab and move the selection ------------
$wsAsis = $wb.Worksheets.Item("asis")
$wsAsis.Activate()
$wsAsis.Range("B11").Select()
